# Update the "as_of_utc" timestamp column (AA) on the "Главные" and
# "Линейные" sheets from 2025-12-04 23:28:42 to 2025-12-05 03:07:02.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-12-05 03:07:02"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Data rows run from row 2 through row 26 (AA1 holds the header).
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA$row").Value = $newTimestamp
    }
}
